$d = $word.ActiveDocument

# The "Towns" mini-table (nested inside the first table) lists town names,
# one of which ("Stewart") must become "Steward". The word is split across
# two runs in the original markup ("Stewar" + "t"); we only want to change
# the trailing "t" to "d", so locate the cell precisely and edit just that
# single trailing character instead of doing a blind document-wide replace
# (which would also rewrite/merge unrelated runs).

$tbl = $d.Tables.Item(1)

foreach ($cell in $tbl.Range.Cells) {
    $txt = $cell.Range.Text
    if ($txt.Length -ge 7 -and $txt.Substring(0, 7) -eq "Stewart") {
        $cellRange = $cell.Range
        $lastChar = $d.Range($cellRange.Start + 6, $cellRange.Start + 7)
        $lastChar.Text = "d"
    }
}
